$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 649.8333
$ws.Range("I100").Value = 649.8333
$ws.Range("K100").Value = 649.8333
$ws.Range("M100").Value = -108.8333
$ws.Range("H113").Value = 4771.143
$ws.Range("I113").Value = 4583.3335
$ws.Range("K113").Value = 4583.3335
$ws.Range("M113").Value = -1329.3335
$ws.Range("H125").Value = 1000
$ws.Range("J125").Value = 1000
$ws.Range("L125").Value = 9000
$ws.Range("N125").Value = -13920
$ws.Range("H132").Value = 3269.102
$ws.Range("I132").Value = 3319.5
$ws.Range("K132").Value = 9958.5
$ws.Range("M132").Value = -7428.5
$ws.Range("H137").Value = 1195.3334
$ws.Range("I137").Value = 1199
$ws.Range("J137").Value = 1193.5
$ws.Range("K137").Value = 3597
$ws.Range("L137").Value = 3580.5
$ws.Range("M137").Value = -1047
$ws.Range("N137").Value = -8680.5
$ws.Range("H138").Value = 3622.1135
$ws.Range("I138").Value = 2149.4119
$ws.Range("J138").Value = 4549.3706
$ws.Range("K138").Value = 6448.2357
$ws.Range("L138").Value = 13648.1118
$ws.Range("M138").Value = -1308.2357
$ws.Range("N138").Value = -23928.1118

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4440.357
$ws.Range("I45").Value = 2777.5
$ws.Range("K45").Value = 2777.5
$ws.Range("M45").Value = -2400.5
$ws.Range("H74").Value = 2724.4736
$ws.Range("J74").Value = 4254.8335
$ws.Range("L74").Value = 4254.8335
$ws.Range("N74").Value = -6002.8335
$ws.Range("H77").Value = 2724.4736
$ws.Range("J77").Value = 4254.8335
$ws.Range("L77").Value = 21274.1675
$ws.Range("N77").Value = -30010.1675
$ws.Range("H132").Value = 2618.7144
$ws.Range("I132").Value = 2656.2424
$ws.Range("J132").Value = 1999.5
$ws.Range("K132").Value = 7968.7272
$ws.Range("L132").Value = 5998.5
$ws.Range("M132").Value = -5438.7272
$ws.Range("N132").Value = -11058.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3553.9565
$ws.Range("I99").Value = 3415.3684
$ws.Range("J99").Value = 4212.25
$ws.Range("K99").Value = 3415.3684
$ws.Range("L99").Value = 4212.25
$ws.Range("M99").Value = -1917.3684
$ws.Range("N99").Value = -7208.25
$ws.Range("H133").Value = 19999
$ws.Range("J133").Value = 19999
$ws.Range("L133").Value = 19999
$ws.Range("N133").Value = -30119
$ws.Range("H134").Value = 2589.25
$ws.Range("I134").Value = 2589.25
$ws.Range("K134").Value = 7767.75
$ws.Range("M134").Value = -5232.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H60").Value = 31945
$ws.Range("J60").Value = 35000
$ws.Range("L60").Value = 35000
$ws.Range("N60").Value = -36022
$ws.Range("H81").Value = 40542.633
$ws.Range("J81").Value = 40542.633
$ws.Range("L81").Value = 40542.633
$ws.Range("N81").Value = -42538.633
$ws.Range("H84").Value = 40542.633
$ws.Range("J84").Value = 40542.633
$ws.Range("L84").Value = 121627.899
$ws.Range("N84").Value = -131611.899
$ws.Range("H98").Value = 70900
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()
$ws.Range("H122").Value = 3344.75
$ws.Range("I122").Value = 3344.75
$ws.Range("K122").Value = 10034.25
$ws.Range("M122").Value = -7584.25
$ws.Range("H125").Value = 89000
$ws.Range("I125").Value = 89000
$ws.Range("K125").Value = 89000
$ws.Range("M125").Value = -86540
$ws.Range("H134").Value = 4166.4165
$ws.Range("I134").Value = 4190
$ws.Range("K134").Value = 12570
$ws.Range("M134").Value = -10035

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H29").Value = 5820.2
$ws.Range("J29").Value = 5820.2
$ws.Range("L29").Value = 17460.6
$ws.Range("N29").Value = -18014.6
$ws.Range("H50").Value = 605
$ws.Range("I50").Value = 490
$ws.Range("J50").Value = 988.3333
$ws.Range("K50").Value = 1470
$ws.Range("L50").Value = 2964.9999
$ws.Range("M50").Value = -989
$ws.Range("N50").Value = -3926.9999
$ws.Range("H53").Value = 605
$ws.Range("I53").Value = 490
$ws.Range("J53").Value = 988.3333
$ws.Range("K53").Value = 1470
$ws.Range("L53").Value = 2964.9999
$ws.Range("M53").Value = -989
$ws.Range("N53").Value = -3926.9999
$ws.Range("H98").Value = 1458
$ws.Range("J98").Value = 1639.6666
$ws.Range("L98").Value = 4918.9998
$ws.Range("N98").Value = -7914.9998
$ws.Range("H121").Value = 683.3333
$ws.Range("J121").Value = 1248.25
$ws.Range("L121").Value = 3744.75
$ws.Range("N121").Value = -6364.75
$ws.Range("H140").Value = 2612
$ws.Range("J140").Value = 3360
$ws.Range("L140").Value = 10080
$ws.Range("N140").Value = -20440

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 28676
$ws.Range("I47").Value = 20028
$ws.Range("K47").Value = 20028
$ws.Range("M47").Value = -19460
$ws.Range("H126").Value = 3617.9375
$ws.Range("I126").Value = 2934.7856
$ws.Range("K126").Value = 8804.356800000001
$ws.Range("M126").Value = -6334.356800000001
$ws.Range("H132").Value = 3963.158
$ws.Range("I132").Value = 4082.625
$ws.Range("J132").Value = 3326
$ws.Range("K132").Value = 12247.875
$ws.Range("L132").Value = 9978
$ws.Range("M132").Value = -9717.875
$ws.Range("N132").Value = -15038

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 8649.825999999999
$ws.Range("I40").Value = 4365.067
$ws.Range("J40").Value = 16683.75
$ws.Range("K40").Value = 4365.067
$ws.Range("L40").Value = 16683.75
$ws.Range("M40").Value = -4229.067
$ws.Range("N40").Value = -16955.75
$ws.Range("H95").Value = 31946.25
$ws.Range("J95").Value = 31946.25
$ws.Range("L95").Value = 31946.25
$ws.Range("N95").Value = -37438.25
$ws.Range("H134").Value = 34999.668
$ws.Range("J134").Value = 34999.668
$ws.Range("L134").Value = 34999.668
$ws.Range("N134").Value = -45139.668
$ws.Range("H136").Value = 1842.7142
$ws.Range("I136").Value = 1779.8
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 5339.4
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -2789.4
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 23924.455
$ws.Range("I81").Value = 5789.3335
$ws.Range("J81").Value = 62785.43
$ws.Range("K81").Value = 11578.667
$ws.Range("L81").Value = 125570.86
$ws.Range("M81").Value = -10517.667
$ws.Range("N81").Value = -127692.86
$ws.Range("H84").Value = 23924.455
$ws.Range("I84").Value = 5789.3335
$ws.Range("J84").Value = 62785.43
$ws.Range("K84").Value = 57893.335
$ws.Range("L84").Value = 627854.3
$ws.Range("M84").Value = -52589.335
$ws.Range("N84").Value = -638462.3
$ws.Range("H96").Value = 1883
$ws.Range("I96").Value = 1879.6
$ws.Range("K96").Value = 1879.6
$ws.Range("M96").Value = -506.5999999999999
$ws.Range("H97").Value = 50266
$ws.Range("J97").Value = 50266
$ws.Range("L97").Value = 50266
$ws.Range("N97").Value = -52248
$ws.Range("H132").Value = 2959.2
$ws.Range("I132").Value = 2998.0527
$ws.Range("J132").Value = 2221
$ws.Range("K132").Value = 8994.158100000001
$ws.Range("L132").Value = 6663
$ws.Range("M132").Value = -6464.158100000001
$ws.Range("N132").Value = -11723
$ws.Range("H136").Value = 19526.092
$ws.Range("I136").Value = 2112.5715
$ws.Range("J136").Value = 49999.75
$ws.Range("K136").Value = 6337.7145
$ws.Range("L136").Value = 149999.25
$ws.Range("M136").Value = -3787.7145
$ws.Range("N136").Value = -155099.25
